# Generate Report for Handback
#
# The localization handback finished: the overall status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", each locale
# sheet records the target file that came back plus the handback xliff
# file name, and the handback timestamp is stamped per locale.

$wb = $excel.ActiveWorkbook

$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/336d9ea627f0cc262fcc8f966ebe33ce27c83d23/e2e/2b4f149c-ea37-4f1f-ac12-f9cd1fd61d1f.md"
$sourceMdName = "2b4f149c-ea37-4f1f-ac12-f9cd1fd61d1f.md"

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# Overview sheet: the status columns for both locales flip to the
# "handed back" status. Widen the two status columns so the longer
# text is readable.
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns("E:E").ColumnWidth = 29.9777047293527
$overview.Columns("F:F").ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------
# zh-cn locale report
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("I2").Value = $sourceMdName
$zhcn.Range("J2").Value = "2b4f149c-ea37-4f1f-ac12-f9cd1fd61d1f.715ac7f8c2e9232e81c2b3ee19f6ad47a4261152.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-27 20:57:40"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceMdUrl, "", "", $sourceMdName)
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = 15570276

$zhcn.Columns("C:C").ColumnWidth = 29.9777047293527
$zhcn.Columns("I:I").ColumnWidth = 39.3234034946987
$zhcn.Columns("J:J").ColumnWidth = 40

# ---------------------------------------------------------------
# de-de locale report
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("I2").Value = $sourceMdName
$dede.Range("J2").Value = "2b4f149c-ea37-4f1f-ac12-f9cd1fd61d1f.715ac7f8c2e9232e81c2b3ee19f6ad47a4261152.de-de.xlf"
$dede.Range("K2").Value = "2016-08-27 20:57:46"

$dede.Hyperlinks.Add($dede.Range("I2"), $sourceMdUrl, "", "", $sourceMdName)
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = 15570276

$dede.Columns("C:C").ColumnWidth = 29.9777047293527
$dede.Columns("I:I").ColumnWidth = 39.3234034946987
$dede.Columns("J:J").ColumnWidth = 40

Write-Output "Handback report generated"
